$wb = $excel.ActiveWorkbook

# ALC row 13
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 3000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 3000
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -3338

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2493.3333
$ws.Range("I43").Value = 1980
$ws.Range("K43").Value = 1980
$ws.Range("M43").Value = -1911

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2450
$ws.Range("I113").Value = 2950
$ws.Range("J113").Value = 2283.3333
$ws.Range("K113").Value = 2950
$ws.Range("L113").Value = 2283.3333
$ws.Range("M113").Value = 304
$ws.Range("N113").Value = -8791.3333

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2258.9429
$ws.Range("I138").Value = 2087.0488
$ws.Range("J138").Value = 2501.9656
$ws.Range("K138").Value = 6261.1464
$ws.Range("L138").Value = 7505.8968
$ws.Range("M138").Value = -1121.1464
$ws.Range("N138").Value = -17785.8968

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1388.7333
$ws.Range("I45").Value = 911.5
$ws.Range("J45").Value = 1706.8889
$ws.Range("K45").Value = 911.5
$ws.Range("L45").Value = 1706.8889
$ws.Range("M45").Value = -534.5
$ws.Range("N45").Value = -2460.8889

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1687.8667
$ws.Range("I122").Value = 1662.2858
$ws.Range("J122").Value = 1710.25
$ws.Range("K122").Value = 4986.857400000001
$ws.Range("L122").Value = 5130.75
$ws.Range("M122").Value = -2536.857400000001
$ws.Range("N122").Value = -10030.75

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1418.4286
$ws.Range("I86").Value = 1361.7778
$ws.Range("J86").Value = 1520.4
$ws.Range("K86").Value = 1361.7778
$ws.Range("L86").Value = 1520.4
$ws.Range("M86").Value = -238.7778000000001
$ws.Range("N86").Value = -3766.4

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1418.4286
$ws.Range("I89").Value = 1361.7778
$ws.Range("J89").Value = 1520.4
$ws.Range("K89").Value = 6808.889
$ws.Range("L89").Value = 7602
$ws.Range("M89").Value = -1192.889
$ws.Range("N89").Value = -18834

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 246.66667
$ws.Range("I94").Value = 246.66667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 246.66667
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 204.33333
$ws.Range("N94").ClearContents()

# BSM row 132
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 70389.5
$ws.Range("J132").Value = 70389.5
$ws.Range("L132").Value = 70389.5
$ws.Range("N132").Value = -80509.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 11195.637
$ws.Range("I134").Value = 16337.714
$ws.Range("J134").Value = 2197
$ws.Range("K134").Value = 49013.142
$ws.Range("L134").Value = 6591
$ws.Range("M134").Value = -46478.142
$ws.Range("N134").Value = -11661

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1488.375
$ws.Range("I10").Value = 1488.375
$ws.Range("K10").Value = 1488.375
$ws.Range("M10").Value = -1349.375

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 6256286.5
$ws.Range("I122").Value = 7360049
$ws.Range("J122").Value = 1633.3334
$ws.Range("K122").Value = 22080147
$ws.Range("L122").Value = 4900.0002
$ws.Range("M122").Value = -22077697
$ws.Range("N122").Value = -9800.0002

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 837.78687
$ws.Range("I5").Value = 543.1429000000001
$ws.Range("J5").Value = 992.475
$ws.Range("K5").Value = 1629.4287
$ws.Range("L5").Value = 2977.425
$ws.Range("M5").Value = -1517.4287
$ws.Range("N5").Value = -3201.425

# CUL row 40
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 137.7
$ws.Range("I40").Value = 130.77777
$ws.Range("K40").Value = 523.11108
$ws.Range("M40").Value = -454.11108

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 217.77777
$ws.Range("I98").Value = 205
$ws.Range("J98").Value = 243.33333
$ws.Range("K98").Value = 615
$ws.Range("L98").Value = 729.99999
$ws.Range("M98").Value = 883
$ws.Range("N98").Value = -3725.99999

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 572.08887
$ws.Range("I113").Value = 589.44446
$ws.Range("J113").Value = 546.05554
$ws.Range("K113").Value = 1768.33338
$ws.Range("L113").Value = 1638.16662
$ws.Range("M113").Value = 401.66662
$ws.Range("N113").Value = -5978.16662

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9804658
$ws.Range("I122").Value = 13333818
$ws.Range("K122").Value = 120004362
$ws.Range("M122").Value = -120001912

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2905.83
$ws.Range("J131").Value = 1871.58
$ws.Range("L131").Value = 5614.74
$ws.Range("N131").Value = -15694.74

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 34483708
$ws.Range("I132").Value = 52632376
$ws.Range("J132").Value = 1237.3
$ws.Range("K132").Value = 473691384
$ws.Range("L132").Value = 11135.7
$ws.Range("M132").Value = -473688854
$ws.Range("N132").Value = -16195.7

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 837.78687
$ws.Range("I135").Value = 543.1429000000001
$ws.Range("J135").Value = 992.475
$ws.Range("K135").Value = 4888.2861
$ws.Range("L135").Value = 8932.275
$ws.Range("M135").Value = -2353.2861
$ws.Range("N135").Value = -14002.275

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3012.875
$ws.Range("I122").Value = 3000.6
$ws.Range("J122").Value = 3033.3333
$ws.Range("K122").Value = 9001.799999999999
$ws.Range("L122").Value = 9099.999899999999
$ws.Range("M122").Value = -6551.799999999999
$ws.Range("N122").Value = -13999.9999

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2001.5834
$ws.Range("I132").Value = 1353.6
$ws.Range("J132").Value = 2464.4285
$ws.Range("K132").Value = 4060.8
$ws.Range("L132").Value = 7393.2855
$ws.Range("M132").Value = -1530.8
$ws.Range("N132").Value = -12453.2855

# LTW row 14
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 892.3077
$ws.Range("I46").Value = 872.7273
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 872.7273
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -684.7273
$ws.Range("N46").Value = -1376

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2669.8
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 2837.5
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 8512.5
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -13412.5

# WVR row 9
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 87505.25
$ws.Range("I9").Value = 50000
$ws.Range("K9").Value = 50000
$ws.Range("M9").Value = -49860

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2973.0557
$ws.Range("I122").Value = 2109.625
$ws.Range("K122").Value = 6328.875
$ws.Range("M122").Value = -3878.875
